# Update classification results (Prediction / Error / Cross Entropy Loss / Success %)
# for the toy NCDE/NODE dataset sheet, per the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9997988619995107
$ws.Range("E2").Value = 0.9997988619995107

# Row 3
$ws.Range("D3").Value = 0.001699055521935858
$ws.Range("E3").Value = 0.001699055521935858

# Row 4
$ws.Range("D4").Value = [double]"2.978673364050297E-07"
$ws.Range("E4").Value = [double]"2.978673364050297E-07"

# Row 5
$ws.Range("D5").Value = 0.0597451618803951
$ws.Range("E5").Value = 0.0597451618803951

# Row 6
$ws.Range("D6").Value = 0.9528007225857484
$ws.Range("E6").Value = 0.9528007225857484

# Row 7
$ws.Range("D7").Value = 0.9999099414284685
$ws.Range("E7").Value = [double]"9.005857153154828E-05"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.2450067696021565
$ws.Range("E8").Value = 0.7549932303978435

# Row 9
$ws.Range("D9").Value = 0.9171529180856638
$ws.Range("E9").Value = 0.08284708191433621

# Row 10
$ws.Range("D10").Value = 0.9977735552026767
$ws.Range("E10").Value = 0.002226444797323257

# Row 11
$ws.Range("D11").Value = 0.8345891892223154
$ws.Range("E11").Value = 0.1654108107776846
$ws.Range("F11").Value = 1.330428719520569
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.999991365074016
$ws.Range("E12").Value = 0.999991365074016

# Row 13
$ws.Range("D13").Value = 0.001113011698155832
$ws.Range("E13").Value = 0.001113011698155832

# Row 14
$ws.Range("D14").Value = [double]"1.368304174964428E-10"
$ws.Range("E14").Value = [double]"1.368304174964428E-10"

# Row 15
$ws.Range("D15").Value = 0.01232078367538299
$ws.Range("E15").Value = 0.01232078367538299

# Row 16
$ws.Range("D16").Value = 0.9961320369616914
$ws.Range("E16").Value = 0.9961320369616914

# Row 17
$ws.Range("D17").Value = 0.9999963210302615
$ws.Range("E17").Value = [double]"3.678969738518667E-06"

# Row 18
$ws.Range("D18").Value = 0.8639090518685172
$ws.Range("E18").Value = 0.1360909481314828

# Row 19
$ws.Range("D19").Value = 0.9906433753326737
$ws.Range("E19").Value = 0.009356624667326341

# Row 20
$ws.Range("D20").Value = 0.9985805355512111
$ws.Range("E20").Value = 0.001419464448788932

# Row 21
$ws.Range("D21").Value = 0.965986732699921
$ws.Range("E21").Value = 0.03401326730007903
$ws.Range("F21").Value = 1.741995215415955
